# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" - F2, F4, F6, F7, F9, F10
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 576
$wsExpo.Range("F4").Value = 392
$wsExpo.Range("F6").Value = 271
$wsExpo.Range("F7").Value = 2447
$wsExpo.Range("F9").Value = 6418
$wsExpo.Range("F10").Value = 172

# Sheet "全部类型" - F2, F4, F6, F9, F11, F12
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 576
$wsAll.Range("F4").Value = 392
$wsAll.Range("F6").Value = 271
$wsAll.Range("F9").Value = 2447
$wsAll.Range("F11").Value = 6418
$wsAll.Range("F12").Value = 172
